# Generate Report for Archive
#
# The localization status report is regenerated: the row for
# "d2739bc9-3c9a-4f06-bc36-c8733fb1ef0d" now sorts ahead of the row for
# "10598507-7115-4d94-9882-232ac59c87d2" (row 5 and row 6 swap which file
# they describe), and the d2739bc9 entry's status has moved back from
# "Ready for handoff" to "In Translation". This is applied identically on
# the Overview sheet and on each per-locale detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (columns: A=File Name, B=zh-cn status, C=de-de status)
$ws = $wb.Worksheets.Item("Overview")

$ws.Cells.Item(5, 1).Value2 = "d2739bc9-3c9a-4f06-bc36-c8733fb1ef0d.md"
$ws.Cells.Item(5, 2).Value2 = "In Translation"
$ws.Cells.Item(5, 3).Value2 = "In Translation"

$ws.Cells.Item(6, 1).Value2 = "10598507-7115-4d94-9882-232ac59c87d2.md"
$ws.Cells.Item(6, 2).Value2 = "Ready for handoff"
$ws.Cells.Item(6, 3).Value2 = "Ready for handoff"

$ws.Hyperlinks.Item(4).TextToDisplay = "d2739bc9-3c9a-4f06-bc36-c8733fb1ef0d.md"
$ws.Hyperlinks.Item(5).TextToDisplay = "10598507-7115-4d94-9882-232ac59c87d2.md"

# ---- Per-locale detail sheets
# columns: A=Source File Name, B=Status, C=Latest Handoff File,
#          D=Latest Handoff Datetime
$locales = @(
    @{ Name = "zh-cn"; HandoffD2739 = "d2739bc9-3c9a-4f06-bc36-c8733fb1ef0d.e50d3ec8988deb4c9204ca2441c33511abf686ef.zh-cn.xlf"; DateD2739 = "2016-03-09 16:30:43"; Handoff10598 = "10598507-7115-4d94-9882-232ac59c87d2.9b6859aac78c6232f04f49ca761b43a6cd26b93d.zh-cn.xlf"; Date10598 = "2016-03-09 16:31:33" },
    @{ Name = "de-de"; HandoffD2739 = "d2739bc9-3c9a-4f06-bc36-c8733fb1ef0d.e50d3ec8988deb4c9204ca2441c33511abf686ef.de-de.xlf"; DateD2739 = "2016-03-09 16:30:47"; Handoff10598 = "10598507-7115-4d94-9882-232ac59c87d2.9b6859aac78c6232f04f49ca761b43a6cd26b93d.de-de.xlf"; Date10598 = "2016-03-09 16:31:43" }
)

foreach ($locale in $locales) {
    $lws = $wb.Worksheets.Item($locale.Name)

    $lws.Cells.Item(5, 1).Value2 = "d2739bc9-3c9a-4f06-bc36-c8733fb1ef0d.md"
    $lws.Cells.Item(5, 2).Value2 = "In Translation"
    $lws.Cells.Item(5, 3).Value2 = $locale.HandoffD2739
    $lws.Cells.Item(5, 4).Value2 = $locale.DateD2739

    $lws.Cells.Item(6, 1).Value2 = "10598507-7115-4d94-9882-232ac59c87d2.md"
    $lws.Cells.Item(6, 2).Value2 = "Ready for handoff"
    $lws.Cells.Item(6, 3).Value2 = $locale.Handoff10598
    $lws.Cells.Item(6, 4).Value2 = $locale.Date10598

    $lws.Hyperlinks.Item(9).TextToDisplay = "d2739bc9-3c9a-4f06-bc36-c8733fb1ef0d.md"
    $lws.Hyperlinks.Item(10).TextToDisplay = $locale.HandoffD2739
    $lws.Hyperlinks.Item(11).TextToDisplay = "10598507-7115-4d94-9882-232ac59c87d2.md"
    $lws.Hyperlinks.Item(12).TextToDisplay = $locale.Handoff10598
}
